# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt - Ciboulette".
# This shifts the existing rows 79..124 down to 80..125 (Excel's native Rows.Insert
# behaviour, copying formatting down with the data) and populates the freshly
# inserted row 79 with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 79, pushing rows 79-124 to 80-125.
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new record.
$ws.Cells.Item(79, 1).Value = 4
$ws.Cells.Item(79, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(79, 3).Value = "Los Lagos"
$ws.Cells.Item(79, 4).Value2 = 44488
$ws.Cells.Item(79, 5).Value = 10
$ws.Cells.Item(79, 6).Value = 100112039
$ws.Cells.Item(79, 7).Value = "Ciboulette"
$ws.Cells.Item(79, 8).Value = "Sin especificar"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 240
$ws.Cells.Item(79, 11).Value = 2500
$ws.Cells.Item(79, 12).Value = 2500
$ws.Cells.Item(79, 13).Value = 2500
$ws.Cells.Item(79, 14).Value = "`$/docena de atados"
$ws.Cells.Item(79, 15).Value = "Región Metropolitana"
$ws.Cells.Item(79, 16).Value = 833
$ws.Cells.Item(79, 17).Value = 3
$ws.Cells.Item(79, 18).Value = "Hortaliza"
